$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Arwa"
$ws.Range("B2").Value = "2023-12-17 17:24:59"
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = "Toqa"
$ws.Range("B3").Value = "2023-12-17 17:25:04"
$ws.Range("C3").Value = 3
